$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9082.916999999999
$ws.Range("I62").Value = 6124.375
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 6124.375
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -5500.375
$ws.Range("N62").Value = -16248

$ws.Range("H64").Value = 4392.75
$ws.Range("I64").Value = 2596.5
$ws.Range("J64").Value = 5470.5
$ws.Range("K64").Value = 2596.5
$ws.Range("L64").Value = 5470.5
$ws.Range("M64").Value = -2348.5
$ws.Range("N64").Value = -5966.5

$ws.Range("H65").Value = 9082.916999999999
$ws.Range("I65").Value = 6124.375
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 30621.875
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -27501.875
$ws.Range("N65").Value = -81240

$ws.Range("H67").Value = 4392.75
$ws.Range("I67").Value = 2596.5
$ws.Range("J67").Value = 5470.5
$ws.Range("K67").Value = 2596.5
$ws.Range("L67").Value = 5470.5
$ws.Range("M67").Value = -1738.5
$ws.Range("N67").Value = -7186.5

$ws.Range("H70").Value = 1561.05
$ws.Range("I70").Value = 1249
$ws.Range("J70").Value = 1694.7858
$ws.Range("K70").Value = 3747
$ws.Range("L70").Value = 5084.357400000001
$ws.Range("M70").Value = -3477
$ws.Range("N70").Value = -5624.357400000001

$ws.Range("H73").Value = 1561.05
$ws.Range("I73").Value = 1249
$ws.Range("J73").Value = 1694.7858
$ws.Range("K73").Value = 3747
$ws.Range("L73").Value = 5084.357400000001
$ws.Range("M73").Value = -2811
$ws.Range("N73").Value = -6956.357400000001

$ws.Range("H124").Value = 41250
$ws.Range("J124").Value = 41250
$ws.Range("L124").Value = 41250
$ws.Range("N124").Value = -51070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3479.2222
$ws.Range("I45").Value = 1500
$ws.Range("J45").Value = 4468.8335
$ws.Range("K45").Value = 1500
$ws.Range("L45").Value = 4468.8335
$ws.Range("M45").Value = -1123
$ws.Range("N45").Value = -5222.8335

$ws.Range("H132").Value = 3814.95
$ws.Range("I132").Value = 3286.889
$ws.Range("J132").Value = 4247
$ws.Range("K132").Value = 9860.667000000001
$ws.Range("L132").Value = 12741
$ws.Range("M132").Value = -7330.667000000001
$ws.Range("N132").Value = -17801

$ws.Range("H133").Value = 43895.4
$ws.Range("J133").Value = 43895.4
$ws.Range("L133").Value = 43895.4
$ws.Range("N133").Value = -48955.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 18100
$ws.Range("I24").Value = 18100
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 18100
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -17865

$ws.Range("H86").Value = 16887.7
$ws.Range("I86").Value = 10559.333
$ws.Range("J86").Value = 26380.25
$ws.Range("K86").Value = 10559.333
$ws.Range("L86").Value = 26380.25
$ws.Range("M86").Value = -9436.333000000001
$ws.Range("N86").Value = -28626.25

$ws.Range("H89").Value = 16887.7
$ws.Range("I89").Value = 10559.333
$ws.Range("J89").Value = 26380.25
$ws.Range("K89").Value = 52796.665
$ws.Range("L89").Value = 131901.25
$ws.Range("M89").Value = -47180.665
$ws.Range("N89").Value = -143133.25

$ws.Range("H133").Value = 34890
$ws.Range("J133").Value = 39853.332
$ws.Range("L133").Value = 39853.332
$ws.Range("N133").Value = -49973.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 1708.1818
$ws.Range("I63").Value = 1618
$ws.Range("J63").Value = 1783.3334
$ws.Range("K63").Value = 4854
$ws.Range("L63").Value = 5350.0002
$ws.Range("M63").Value = -4105
$ws.Range("N63").Value = -6848.0002

$ws.Range("H66").Value = 1708.1818
$ws.Range("I66").Value = 1618
$ws.Range("J66").Value = 1783.3334
$ws.Range("K66").Value = 14562
$ws.Range("L66").Value = 16050.0006
$ws.Range("M66").Value = -10818
$ws.Range("N66").Value = -23538.0006

$ws.Range("H68").Value = 2814.5522
$ws.Range("I68").Value = 3757.7317
$ws.Range("J68").Value = 1327.2307
$ws.Range("K68").Value = 11273.1951
$ws.Range("L68").Value = 3981.6921
$ws.Range("M68").Value = -10462.1951
$ws.Range("N68").Value = -5603.6921

$ws.Range("H71").Value = 2814.5522
$ws.Range("I71").Value = 3757.7317
$ws.Range("J71").Value = 1327.2307
$ws.Range("K71").Value = 33819.5853
$ws.Range("L71").Value = 11945.0763
$ws.Range("M71").Value = -29763.5853
$ws.Range("N71").Value = -20057.0763

$ws.Range("H114").Value = 1250.2941
$ws.Range("I114").Value = 469.75
$ws.Range("J114").Value = 1944.1111
$ws.Range("K114").Value = 1409.25
$ws.Range("L114").Value = 5832.3333
$ws.Range("M114").Value = 1844.75
$ws.Range("N114").Value = -12340.3333

$ws.Range("H116").Value = 1332.8572
$ws.Range("I116").Value = 959.6
$ws.Range("K116").Value = 2878.8
$ws.Range("M116").Value = 563.1999999999998

$ws.Range("H131").Value = 989.6
$ws.Range("I131").Value = 427.36365
$ws.Range("J131").Value = 1059.0898
$ws.Range("K131").Value = 1282.09095
$ws.Range("L131").Value = 3177.2694
$ws.Range("M131").Value = 3757.90905
$ws.Range("N131").Value = -13257.2694

$ws.Range("H132").Value = 764.3333
$ws.Range("I132").Value = 577.6
$ws.Range("J132").Value = 997.75
$ws.Range("K132").Value = 5198.400000000001
$ws.Range("L132").Value = 8979.75
$ws.Range("M132").Value = -2668.400000000001
$ws.Range("N132").Value = -14039.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

$ws.Range("H102").Value = 2428.5
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 2857
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 2857
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6101

$ws.Range("H122").Value = 2778654.5
$ws.Range("I122").Value = 2778654.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8335963.5
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -8333513.5

$ws.Range("H132").Value = 2737.7437
$ws.Range("I132").Value = 2259.3103
$ws.Range("J132").Value = 4125.2
$ws.Range("K132").Value = 6777.9309
$ws.Range("L132").Value = 12375.6
$ws.Range("M132").Value = -4247.9309
$ws.Range("N132").Value = -17435.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8251
$ws.Range("I61").Value = 10001.333
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 10001.333
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -9799.333000000001
$ws.Range("N61").Value = -3404

$ws.Range("H113").Value = 8251
$ws.Range("I113").Value = 10001.333
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 10001.333
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -7831.333000000001
$ws.Range("N113").Value = -7340

$ws.Range("H122").Value = 3593.6875
$ws.Range("I122").Value = 1833
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5499
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3049
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 3163.2144
$ws.Range("I132").Value = 2416.1738
$ws.Range("K132").Value = 7248.5214
$ws.Range("M132").Value = -4718.5214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 100000000
$ws.Range("I2").Value = 100000000
$ws.Range("K2").Value = 100000000
$ws.Range("M2").Value = -99999888

$ws.Range("H4").Value = 25801
$ws.Range("I4").Value = 7501
$ws.Range("K4").Value = 7501
$ws.Range("M4").Value = -7388

$ws.Range("H113").Value = 1374.5
$ws.Range("I113").Value = 1456.5714
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 4369.7142
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -2199.7142
$ws.Range("N113").Value = -6740

$ws.Range("H122").Value = 41680.08
$ws.Range("I122").Value = 112265.78
$ws.Range("J122").Value = 1975.625
$ws.Range("K122").Value = 336797.34
$ws.Range("L122").Value = 5926.875
$ws.Range("M122").Value = -334347.34
$ws.Range("N122").Value = -10826.875

$ws.Range("H132").Value = 12822714
$ws.Range("I132").Value = 17243040
$ws.Range("J132").Value = 3765.5
$ws.Range("K132").Value = 51729120
$ws.Range("L132").Value = 11296.5
$ws.Range("M132").Value = -51726590
$ws.Range("N132").Value = -16356.5

$ws.Range("H135").Value = 98233.92
$ws.Range("J135").Value = 98233.92
$ws.Range("L135").Value = 98233.92
$ws.Range("N135").Value = -108373.92
